# Update "want to go" counts (column F) on the 展览 (Exhibition) sheet
# and the corresponding rows on the 全部类型 (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExpo.Range("F3").Value  = 2149
$wsExpo.Range("F5").Value  = 11252
$wsExpo.Range("F10").Value = 11173
$wsExpo.Range("F13").Value = 55
$wsExpo.Range("F15").Value = 5592
$wsExpo.Range("F16").Value = 98
$wsExpo.Range("F17").Value = 3451
$wsExpo.Range("F18").Value = 8

# 全部类型 sheet updates (same events, different rows)
$wsAll.Range("F3").Value  = 2149
$wsAll.Range("F7").Value  = 11252
$wsAll.Range("F12").Value = 11173
$wsAll.Range("F15").Value = 55
$wsAll.Range("F17").Value = 5592
$wsAll.Range("F18").Value = 98
$wsAll.Range("F19").Value = 3451
$wsAll.Range("F20").Value = 8
